$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.523.62'
$ws.Range("E2").Value = '  +5.17%  '
$ws.Range("D3").Value = '3.172.79'
$ws.Range("E3").Value = '  +2.03%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '''400.45'
$ws.Range("E5").Value = '  +2.80%  '
$ws.Range("D6").Value = '''109.64'
$ws.Range("E6").Value = '  +5.59%  '
$ws.Range("D7").Value = '''0.549'
$ws.Range("E7").Value = '  +0.65%  '
$ws.Range("D8").Value = '''1.00'
$ws.Range("D9").Value = '''0.616'
$ws.Range("E9").Value = '  +4.08%  '
$ws.Range("D10").Value = '''39.13'
$ws.Range("E10").Value = '  +5.06%  '
$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").Value = '''0.140'
$ws.Range("E11").Value = '  +1.69%  '
$ws.Range("B12").Value = 'Dogecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D12").Value = '''0.0885'
$ws.Range("E12").Value = '  +2.61%  '
$ws.Range("D13").Value = '3.684.47'
$ws.Range("E13").Value = '  +2.38%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = '''19.05'
$ws.Range("E14").Value = '  +1.58%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = '''8.04'
$ws.Range("E15").Value = '  +1.84%  '
$ws.Range("E16").Value = '  +7.27%  '
$ws.Range("D17").Value = '3.161.06'
$ws.Range("E17").Value = '  +2.17%  '
$ws.Range("E18").Value = '  -2.66%  '
$ws.Range("D19").Value = '54.466.57'
$ws.Range("E19").Value = '  +4.93%  '
$ws.Range("E20").Value = '  +2.69%  '
$ws.Range("E21").Value = '  +3.78%  '
$ws.Range("E22").Value = '  +2.88%  '
$ws.Range("D23").Value = '''72.15'
$ws.Range("E23").Value = '  +3.06%  '
$ws.Range("D24").Value = '''274.84'
$ws.Range("E24").Value = '  +2.36%  '
$ws.Range("D25").Value = '''3.27'
$ws.Range("E25").Value = '  +4.58%  '
$ws.Range("D26").Value = '''7.96'
$ws.Range("E26").Value = '  -2.04%  '
$ws.Range("D27").Value = '''27.83'
$ws.Range("E27").Value = '  +2.45%  '
$ws.Range("D28").Value = '''7.59'
$ws.Range("E28").Value = '  +5.06%  '
$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").Value = '''0.169'
$ws.Range("E29").Value = '  -0.75%  '
$ws.Range("B30").Value = 'Dai'
$ws.Range("C30").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D30").Value = '''0.999'
$ws.Range("E30").Value = '  -0.09%  '
$ws.Range("E31").Value = '  +2.37%  '
$ws.Range("E32").Value = '  +6.33%  '
$ws.Range("D33").Value = '''0.0510'
$ws.Range("E33").Value = '  +12.91%  '
$ws.Range("D34").Value = '''36.59'
$ws.Range("E34").Value = '  +3.16%  '
$ws.Range("E35").Value = '  +1.38%  '
$ws.Range("D36").Value = '''51.53'
$ws.Range("E36").Value = '  +2.15%  '
$ws.Range("D37").Value = '''3.62'
$ws.Range("E37").Value = '  +5.81%  '
$ws.Range("D38").Value = '''1.00'
$ws.Range("E38").Value = '  +0.03%  '
$ws.Range("D39").Value = '''2.90'
$ws.Range("E39").Value = '  +11.02%  '
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").Value = '''1.94'
$ws.Range("E40").Value = '  +2.53%  '
$ws.Range("B41").Value = 'NEARProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D41").Value = '''4.08'
$ws.Range("E41").Value = '  +9.92%  '
$ws.Range("E42").Value = '  +0.97%  '
$ws.Range("D43").Value = '''17.23'
$ws.Range("E43").Value = '  +1.86%  '
$ws.Range("D44").Value = '''131.27'
$ws.Range("E44").Value = '  +1.72%  '
$ws.Range("E45").Value = '  +0.95%  '
$ws.Range("D46").Value = '''22.12'
$ws.Range("E46").Value = '  -0.93%  '
$ws.Range("D47").Value = '''2.46'
$ws.Range("E47").Value = '  -1.94%  '
$ws.Range("E48").Value = '  -0.60%  '
$ws.Range("D49").Value = '2.093.77'
$ws.Range("E49").Value = '  +2.12%  '
$ws.Range("D50").Value = '''0.0347'
$ws.Range("E50").Value = '  +6.48%  '
$ws.Range("D51").Value = '''0.0508'
$ws.Range("E51").Value = '  +12.73%  '
